# Update the RBSP seminar title slide for the next talk:
#  - swap speaker names / talk titles in the "Title 1" placeholder
#    (and bold the two speaker-name runs)
#  - grow the "Title 1" placeholder to fit the extra text
#  - swap the series subtitle text in the "Subtitle 2" placeholder

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---- Title placeholder (speakers + talk titles) -----------------------
$titleShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Title 1") { $titleShape = $sh }
}

$tr = $titleShape.TextFrame.TextRange
$delta = 0

# Run: "Matina Gkioulidou:" -> "Ian Richardson:" (bold)
$start = 1 + $delta
$c = $tr.Characters($start, 18)
$oldLen = $c.Text.Length
$newText = "Ian Richardson:"
$c.Text = $newText
$delta += ($newText.Length - $oldLen)
$tr.Characters($start, $newText.Length).Font.Bold = $true

# Run: "Observing the Global Geospace in Mesoscale Resolution" ->
#      "A Heliospheric View of Solar Energetic Particle Events "
$start = 20 + $delta
$c = $tr.Characters($start, 53)
$oldLen = $c.Text.Length
$newText = "A Heliospheric View of Solar Energetic Particle Events "
$c.Text = $newText
$delta += ($newText.Length - $oldLen)

# Run: "Anthony " -> "Rachael " (bold)
$start = 75 + $delta
$c = $tr.Characters($start, 8)
$oldLen = $c.Text.Length
$newText = "Rachael "
$c.Text = $newText
$delta += ($newText.Length - $oldLen)
$tr.Characters($start, $newText.Length).Font.Bold = $true

# Run: "Sciola" -> "Filwett" (bold, keeps its err="1" spellcheck flag)
$start = 83 + $delta
$c = $tr.Characters($start, 6)
$oldLen = $c.Text.Length
$newText = "Filwett"
$c.Text = $newText
$delta += ($newText.Length - $oldLen)
$tr.Characters($start, $newText.Length).Font.Bold = $true

# Run: ":" stays ":" but becomes bold
$start = 89 + $delta
$tr.Characters($start, 1).Font.Bold = $true

# Run: "Build-up of the Storm-Time Ring Current via Mesoscale Plasma Sheet Flows" ->
#      "Solar Energetic Particle Access to the Magnetosphere, a Comparison of 4 SEP Events Measured with RBSP"
$start = 91 + $delta
$c = $tr.Characters($start, 72)
$oldLen = $c.Text.Length
$newText = "Solar Energetic Particle Access to the Magnetosphere, a Comparison of 4 SEP Events Measured with RBSP"
$c.Text = $newText
$delta += ($newText.Length - $oldLen)

# Grow the title placeholder so the extra text still fits.
$titleShape.Height = 351.6086

# ---- Subtitle placeholder (series title) -------------------------------
$subtitleShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Subtitle 2") { $subtitleShape = $sh }
}

$subTr = $subtitleShape.TextFrame.TextRange
$subC = $subTr.Characters(1, $subTr.Text.Length)
$subC.Text = "Solar Energetic Particles in Heliosphere and Geospace"
